$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-19 Sunday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-05-20 Monday", 2) | Out-Null
$d.Content.Find.Execute("757÷2=378, 1", $true, $true, $false, $false, $false, $true, 1, $false, "820÷3=273, 1", 2) | Out-Null
$d.Content.Find.Execute("255÷3=85, 0", $true, $true, $false, $false, $false, $true, 1, $false, "713÷5=142, 3", 2) | Out-Null
$d.Content.Find.Execute("678÷3=226, 0", $true, $true, $false, $false, $false, $true, 1, $false, "586÷7=83, 5", 2) | Out-Null
$d.Content.Find.Execute("297÷5=59, 2", $true, $true, $false, $false, $false, $true, 1, $false, "826÷9=91, 7", 2) | Out-Null
$d.Content.Find.Execute("493÷5=98, 3", $true, $true, $false, $false, $false, $true, 1, $false, "526÷2=263, 0", 2) | Out-Null
$d.Content.Find.Execute("368÷8=46, 0", $true, $true, $false, $false, $false, $true, 1, $false, "927÷4=231, 3", 2) | Out-Null
$d.Content.Find.Execute("246÷5=49, 1", $true, $true, $false, $false, $false, $true, 1, $false, "444÷8=55, 4", 2) | Out-Null
$d.Content.Find.Execute("835÷4=208, 3", $true, $true, $false, $false, $false, $true, 1, $false, "129÷4=32, 1", 2) | Out-Null
$d.Content.Find.Execute("791÷8=98, 7", $true, $true, $false, $false, $false, $true, 1, $false, "228÷7=32, 4", 2) | Out-Null
$d.Content.Find.Execute("466÷4=116, 2", $true, $true, $false, $false, $false, $true, 1, $false, "806÷2=403, 0", 2) | Out-Null
$d.Content.Find.Execute("670÷4=167, 2", $true, $true, $false, $false, $false, $true, 1, $false, "698÷4=174, 2", 2) | Out-Null
$d.Content.Find.Execute("579÷7=82, 5", $true, $true, $false, $false, $false, $true, 1, $false, "888÷4=222, 0", 2) | Out-Null
$d.Content.Find.Execute("432÷7=61, 5", $true, $true, $false, $false, $false, $true, 1, $false, "688÷6=114, 4", 2) | Out-Null
$d.Content.Find.Execute("274÷6=45, 4", $true, $true, $false, $false, $false, $true, 1, $false, "126÷6=21, 0", 2) | Out-Null
$d.Content.Find.Execute("337÷5=67, 2", $true, $true, $false, $false, $false, $true, 1, $false, "454÷7=64, 6", 2) | Out-Null
$d.Content.Find.Execute("271÷7=38, 5", $true, $true, $false, $false, $false, $true, 1, $false, "663÷5=132, 3", 2) | Out-Null
$d.Content.Find.Execute("888÷7=126, 6", $true, $true, $false, $false, $false, $true, 1, $false, "845÷2=422, 1", 2) | Out-Null
$d.Content.Find.Execute("521÷2=260, 1", $true, $true, $false, $false, $false, $true, 1, $false, "127÷8=15, 7", 2) | Out-Null
$d.Content.Find.Execute("643÷3=214, 1", $true, $true, $false, $false, $false, $true, 1, $false, "249÷5=49, 4", 2) | Out-Null
$d.Content.Find.Execute("785÷9=87, 2", $true, $true, $false, $false, $false, $true, 1, $false, "999÷6=166, 3", 2) | Out-Null
$d.Content.Find.Execute("818÷4=204, 2", $true, $true, $false, $false, $false, $true, 1, $false, "228÷6=38, 0", 2) | Out-Null
$d.Content.Find.Execute("716÷8=89, 4", $true, $true, $false, $false, $false, $true, 1, $false, "836÷5=167, 1", 2) | Out-Null
$d.Content.Find.Execute("218÷5=43, 3", $true, $true, $false, $false, $false, $true, 1, $false, "317÷4=79, 1", 2) | Out-Null
$d.Content.Find.Execute("192÷3=64, 0", $true, $true, $false, $false, $false, $true, 1, $false, "270÷5=54, 0", 2) | Out-Null
$d.Content.Find.Execute("242÷8=30, 2", $true, $true, $false, $false, $false, $true, 1, $false, "185÷2=92, 1", 2) | Out-Null
